# Sprint Backlog slide: update the "Gerenciar produtos" bullet list.
#   - "Consultar disponibilidade" -> "Consultar disponibilidade (por nome ou código)"
#   - "Encomendar produto"        -> "Alterar um produto"
#   - add a new bullet "Excluir produto" right after it

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(6)
$sh = $s.Shapes.Item(1)

# --- Edit 1: "Consultar disponibilidade" -> "Consultar disponibilidade (por nome ou código)"
$tr    = $sh.TextFrame.TextRange
$para5 = $tr.Paragraphs(5, 1)
$start5 = $para5.Start
$len5   = $para5.Length
# Insert the full replacement text right after the old paragraph text (this copies the
# existing run's formatting/rPr), then remove the old text, leaving a single clean run.
$para5.InsertAfter("Consultar disponibilidade (por nome ou código)") | Out-Null
$sh.TextFrame.TextRange.Characters($start5, $len5 - 1).Delete()

# --- Edit 2: "Encomendar produto" -> "Alterar um produto"
$tr    = $sh.TextFrame.TextRange
$para6 = $tr.Paragraphs(6, 1)
$start6 = $para6.Start
$len6   = $para6.Length
$para6.InsertAfter("Alterar um produto") | Out-Null
$sh.TextFrame.TextRange.Characters($start6, $len6 - 1).Delete()

# --- Edit 3: add a new bullet paragraph "Excluir produto" right after "Alterar um produto"
$tr     = $sh.TextFrame.TextRange
$para6b = $tr.Paragraphs(6, 1)
$para6b.InsertAfter([char]13 + "Excluir produto") | Out-Null
